$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row (data rows below header)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Append "*" to the PettittBreak (column B) values and store as text,
# before we shift columns around.
for ($r = 2; $r -le $lastRow; $r++) {
    $val = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 2).Value = "$val*"
}

# Remove column C (PettittP) entirely; this shifts column D (TaylorBreak)
# left into column C.
$ws.Range("C:C").Delete(-4159)
